# Weekly NYPD CompStat refresh: new crime data collected.
# Updates the report title (volume/number + date range) and the
# Week/28-Day/YTD/2-Yr/15-Yr/32-Yr crime statistics table (rows 15-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Masthead: volume number and the date range covered by the report.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/1/2025  Through  12/7/2025"

# ---------------------------------------------------------------------
# Helper cells used purely as format/value donors for PasteSpecial so
# that converted cells land on the workbook's existing shared styles
# instead of minting new (duplicate) style entries. None of these
# donor cells are themselves modified by this week's refresh.
#   C14 -> text "0"      (style: plain/General, font 7)
#   E14 -> text "***.*"  (style: plain/General, font 7)
#   I15 -> number style, integer columns  (numFmt #,##0)
#   M15 -> number style, percent columns  (numFmt #,##0.0;"-"#,##0.0)
# ---------------------------------------------------------------------
$fmtText0   = $ws.Range("C14")
$fmtTextNA  = $ws.Range("E14")
$fmtNumInt  = $ws.Range("I15")
$fmtNumPct  = $ws.Range("M15")

function Set-NumericCell($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $fmtNumInt.Copy()
    $cell.PasteSpecial(-4122)
    $cell.Value = $value
}

function Set-PercentCell($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $fmtNumPct.Copy()
    $cell.PasteSpecial(-4122)
    $cell.Value = $value
}

function Set-Text0Cell($cellRef) {
    $cell = $ws.Range($cellRef)
    $fmtText0.Copy()
    $cell.PasteSpecial(-4122)
    $fmtText0.Copy()
    $cell.PasteSpecial(-4163)
}

function Set-TextNACell($cellRef) {
    $cell = $ws.Range($cellRef)
    $fmtTextNA.Copy()
    $cell.PasteSpecial(-4122)
    $fmtTextNA.Copy()
    $cell.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
Set-NumericCell "D15" 4
Set-PercentCell "E15" -100
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -80
$ws.Range("J15").Value = 17
$ws.Range("K15").Value = -5.882352941176
$ws.Range("L15").Value = 166.666666666667

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("F16").Value = 5
$ws.Range("H16").Value = -44.444444444444
$ws.Range("I16").Value = 101
$ws.Range("J16").Value = 105
$ws.Range("K16").Value = -3.809523809523
$ws.Range("L16").Value = -20.472440944881
$ws.Range("M16").Value = -5.607476635514
$ws.Range("N16").Value = -82.525951557093

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 20
$ws.Range("I17").Value = 149
$ws.Range("J17").Value = 125
$ws.Range("K17").Value = 19.2
$ws.Range("L17").Value = -1.324503311258
$ws.Range("M17").Value = 18.253968253968
$ws.Range("N17").Value = -36.324786324786

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
Set-NumericCell "D18" 7
Set-PercentCell "E18" -100
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -91.666666666666
$ws.Range("J18").Value = 94
$ws.Range("K18").Value = -1.063829787234
$ws.Range("L18").Value = -29.545454545454
$ws.Range("M18").Value = -1.063829787234
$ws.Range("N18").Value = -84.703947368421

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -57.142857142857
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = -17.391304347826
$ws.Range("I19").Value = 603
$ws.Range("J19").Value = 660
$ws.Range("K19").Value = -8.636363636363
$ws.Range("L19").Value = -16.133518776077
$ws.Range("M19").Value = -0.658978583196
$ws.Range("N19").Value = -22.890025575447

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
Set-Text0Cell "C20"
Set-Text0Cell "D20"
Set-TextNACell "E20"
$ws.Range("F20").Value = 2
$ws.Range("H20").Value = -50
$ws.Range("L20").Value = -58.666666666666
$ws.Range("M20").Value = -22.5
$ws.Range("N20").Value = -93.528183716075

# ---------------------------------------------------------------------
# Row 21 - TOTAL (bold)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -61.290322580645
$ws.Range("F21").Value = 59
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = -32.183908045977
$ws.Range("I21").Value = 994
$ws.Range("J21").Value = 1051
$ws.Range("K21").Value = -5.423406279733
$ws.Range("L21").Value = -17.851239669421
$ws.Range("M21").Value = 1.016260162601
$ws.Range("N21").Value = -63.361592333210

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
Set-NumericCell "D22" 1
Set-PercentCell "E22" -100
$ws.Range("J22").Value = 37
$ws.Range("K22").Value = -27.027027027027
$ws.Range("L22").Value = -25

# ---------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------
$ws.Range("C23").Value = 2
Set-NumericCell "D23" 4
Set-PercentCell "E23" -50
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 44
$ws.Range("J23").Value = 44
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = -8.333333333333
$ws.Range("M23").Value = 7.317073170731

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = -40
$ws.Range("F24").Value = 63
$ws.Range("G24").Value = 95
$ws.Range("H24").Value = -33.684210526315
$ws.Range("I24").Value = 810
$ws.Range("J24").Value = 724
$ws.Range("K24").Value = 11.878453038674
$ws.Range("L24").Value = 4.381443298969
$ws.Range("M24").Value = 0.746268656716

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 36
$ws.Range("H25").Value = -36.842105263157
$ws.Range("I25").Value = 433
$ws.Range("J25").Value = 360
$ws.Range("K25").Value = 20.277777777777
$ws.Range("L25").Value = 1.405152224824

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 13
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = -43.478260869565
$ws.Range("I26").Value = 282
$ws.Range("J26").Value = 310
$ws.Range("K26").Value = -9.032258064516
$ws.Range("L26").Value = -16.814159292035
$ws.Range("M26").Value = -9.032258064516

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------
Set-NumericCell "D27" 4
Set-PercentCell "E27" -100
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -80
$ws.Range("J27").Value = 28
$ws.Range("K27").Value = -35.714285714285
$ws.Range("L27").Value = 100

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
Set-PercentCell "E28" 100
$ws.Range("F28").Value = 4
Set-NumericCell "G28" 1
Set-PercentCell "H28" 300
$ws.Range("I28").Value = 53
$ws.Range("J28").Value = 51
$ws.Range("K28").Value = 3.921568627450
$ws.Range("L28").Value = -5.357142857142

# ---------------------------------------------------------------------
# Row 31 - Hate Crimes
# ---------------------------------------------------------------------
Set-Text0Cell "G31"
Set-TextNACell "H31"
$ws.Range("L31").Value = -66.666666666666
